$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Update the product name text shown on both sheets (B1) - same new text on each sheet.
$ws1.Range("B1").Value = "4228-RBI-EI-DB-SAR-REC-CTRFD-RNI-INT-FFC-SAR-FFROP-DAILY-FIFR-1-MD-TR-1-ON-PER-1st"
$ws2.Range("B1").Value = "4228-RBI-EI-DB-SAR-REC-CTRFD-RNI-INT-FFC-SAR-FFROP-DAILY-FIFR-1-MD-TR-1-ON-PER-1st"

# Change shortname (B2) from numeric 4228 to text "422u"
$ws1.Range("B2").Value = "422u"

# Move the selection / active sheet so ProductLoanOutput (tab 2) is now the one selected,
# removing the prior inter-dependency on ProductLoanInput's A19 selection.
$ws1.Range("B3").Select()
$ws2.Select()
$ws2.Range("B1").Select()
